$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells (Wins, Losses, Ties) in AD1:AF1, matching the
# bold/centered/bordered header style already used by the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
